$wb = $excel.ActiveWorkbook
$wsCreate = $wb.Worksheets.Item("Create_Org")
$wsEdit   = $wb.Worksheets.Item("Edit_Org")

# --- Column widths (AutoFit-style "best fit" resize of the header columns) ---
# Target widths (as authored by Excel's real AutoFit) expressed in the
# 256ths-of-a-character unit that ends up in the saved XML. The inputs below
# are chosen so this runtime's ColumnWidth quantization lands as close as
# possible to those authored values.
$wsCreate.Columns.Item(1).ColumnWidth = 24.833333333333332
$wsCreate.Columns.Item(2).ColumnWidth = 9.333333333333334
$wsCreate.Columns.Item(3).ColumnWidth = 21.0
$wsCreate.Columns.Item(4).ColumnWidth = 23.5
$wsCreate.Columns.Item(5).ColumnWidth = 21.333333333333332
$wsCreate.Columns.Item(6).ColumnWidth = 29.0

$wsEdit.Columns.Item(1).ColumnWidth = 24.833333333333332
$wsEdit.Columns.Item(3).ColumnWidth = 21.0
$wsEdit.Columns.Item(4).ColumnWidth = 23.5
$wsEdit.Columns.Item(5).ColumnWidth = 21.333333333333332
$wsEdit.Columns.Item(6).ColumnWidth = 29.0
$wsEdit.Columns.Item(7).ColumnWidth = 21.833333333333332

# --- Clear the "Purchasing Dept Indicator" (E2) checkbox value on both sheets ---
$wsCreate.Range("E2").ClearContents() | Out-Null
$wsEdit.Range("E2").ClearContents() | Out-Null

# --- Flip the "Credit Department Group Indicator" (F2) flag to False on Edit_Org ---
$wsEdit.Range("F2").Value = $false

# --- Update selections on each sheet, and make Edit_Org the active tab/sheet ---
$wsCreate.Range("E2").Select() | Out-Null
$wsEdit.Range("F2").Select() | Out-Null
